$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 16 ("only 5 comments is shown...")
# so the sprint1/models plan rows shift down by one, mirroring the Excel
# "insert row" behaviour (copies formatting from the row above).
$ws.Rows.Item(16).Insert()

# Populate the new row's D cell with the new comment plan text.
$ws.Range("D16").Value = "can give comment to the post or specific comment"

# Update the active selection to match the authored change (D7 -> D10).
[void]$ws.Range("D10").Select()
